# Saldo.xlsx update:
#  - Remove the data row for account 005366671 / TATIANA (spreadsheet row 2),
#    shifting all subsequent rows up by one.
#  - Update the Saldo for account 005064129 / THIAGO (now row 2 after the
#    deletion) from 20777.32 to 18791.64.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire TATIANA row (row 2); remaining rows shift up.
$ws.Rows.Item(2).Delete()

# THIAGO's row is now row 2; update the Saldo (column C) value.
$ws.Range("C2").Value = 18791.64
